$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 198 with new values (Actualizacion desde MV -datos-) ---
$ws.Range("B198").Value = 1.48
$ws.Range("C198").Value = 1.01
$ws.Range("I198").Value = 2.1
$ws.Range("J198").Value = 7.45
$ws.Range("L198").Value = 6.78
$ws.Range("M198").Value = 11.08
$ws.Range("N198").Value = 7.72
$ws.Range("O198").Value = 7.48
$ws.Range("P198").Value = 6.32

# --- Add new date rows 199-204 as text (shared strings), matching existing "Serie" column format ---
$ws.Range("A199:A204").NumberFormat = "@"
$ws.Range("A199").Value = "05-10-2021"
$ws.Range("A200").Value = "06-10-2021"
$ws.Range("A201").Value = "07-10-2021"
$ws.Range("A202").Value = "08-10-2021"
$ws.Range("A203").Value = "11-10-2021"
$ws.Range("A204").Value = "12-10-2021"
$ws.Range("A199:A204").ClearFormats()

# Row 199
$ws.Range("B199").Value = 1.53
$ws.Range("C199").Value = 1.08
$ws.Range("D199").Value = 0.06
$ws.Range("E199").Value = -0.19
$ws.Range("F199").Value = 2.33
$ws.Range("G199").Value = 3.44
$ws.Range("H199").Value = 2.25
$ws.Range("I199").Value = 2.1
$ws.Range("J199").Value = 7.55
$ws.Range("K199").Value = 1.88
$ws.Range("L199").Value = 6.8
$ws.Range("M199").Value = 11.14
$ws.Range("O199").Value = 7.56
$ws.Range("P199").Value = 6.28

# Row 200
$ws.Range("B200").Value = 1.52
$ws.Range("C200").Value = 1.07
$ws.Range("D200").Value = 0.08
$ws.Range("E200").Value = -0.18
$ws.Range("F200").Value = 2.42
$ws.Range("G200").Value = 3.57
$ws.Range("H200").Value = 2.38
$ws.Range("I200").Value = 2.09
$ws.Range("J200").Value = 7.58
$ws.Range("K200").Value = 1.89
$ws.Range("L200").Value = 6.84
$ws.Range("M200").Value = 11.04
$ws.Range("O200").Value = 7.54
$ws.Range("P200").Value = 6.12

# Row 201
$ws.Range("B201").Value = 1.57
$ws.Range("C201").Value = 1.08
$ws.Range("D201").Value = 0.07000000000000001
$ws.Range("E201").Value = -0.19
$ws.Range("F201").Value = 2.4
$ws.Range("G201").Value = 3.58
$ws.Range("H201").Value = 2.41
$ws.Range("I201").Value = 2.08
$ws.Range("J201").Value = 7.56
$ws.Range("K201").Value = 1.92
$ws.Range("L201").Value = 6.74
$ws.Range("M201").Value = 11.09
$ws.Range("O201").Value = 7.5
$ws.Range("P201").Value = 5.91

# Row 202
$ws.Range("B202").Value = 1.61
$ws.Range("C202").Value = 1.16
$ws.Range("D202").Value = 0.09
$ws.Range("E202").Value = -0.15
$ws.Range("F202").Value = 2.39
$ws.Range("G202").Value = 3.62
$ws.Range("H202").Value = 2.42
$ws.Range("I202").Value = 2.08
$ws.Range("J202").Value = 7.5
$ws.Range("K202").Value = 1.94
$ws.Range("L202").Value = 6.84
$ws.Range("M202").Value = 10.88
$ws.Range("O202").Value = 7.49

# Row 203
$ws.Range("B203").Value = 1.61
$ws.Range("C203").Value = 1.19
$ws.Range("D203").Value = 0.09
$ws.Range("E203").Value = -0.12
$ws.Range("G203").Value = 3.67
$ws.Range("H203").Value = 2.52
$ws.Range("I203").Value = 2.17
$ws.Range("J203").Value = 7.53
$ws.Range("K203").Value = 1.97
$ws.Range("L203").Value = 6.84
$ws.Range("M203").Value = 11
$ws.Range("O203").Value = 7.5
$ws.Range("P203").Value = 5.96

# Row 204
$ws.Range("B204").Value = 1.61
$ws.Range("C204").Value = 1.17
$ws.Range("D204").Value = 0.09
$ws.Range("E204").Value = -0.12
$ws.Range("F204").Value = 2.45
$ws.Range("G204").Value = 3.61
$ws.Range("H204").Value = 2.63
$ws.Range("I204").Value = 2.27
$ws.Range("J204").Value = 7.55
$ws.Range("K204").Value = 1.99
$ws.Range("L204").Value = 6.84
$ws.Range("O204").Value = 7.49
$ws.Range("P204").Value = 5.97
